# Refresh the crypto price list (Sheet1, rows 2-51) to match the latest scrape.
# Columns D (Price) and E (Volume/1h) are stored as plain text in this workbook
# (values use "." as a thousands separator and keep padded "%" strings), so a
# handful of single-decimal price strings that Excel would otherwise silently
# reinterpret as numbers are forced to Text format before being written and the
# style is put back to Normal immediately after so no formatting residue remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.635.80'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '2.206.84'
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.615'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.75'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.22%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.587'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.87'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0908'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.86'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.102'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').Value = '2.539.80'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.82%  '
$ws.Range('D16').Value = '2.205.56'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.775'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.85%  '
$ws.Range('D18').Value = '42.580.01'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000102'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.03'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('E23').Value = '  -3.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.36'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -8.61%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.54'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.27%  '
$ws.Range('B27').Value = 'WEMIXToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.01%  '
$ws.Range('E29').Value = '  +2.88%  '
$ws.Range('E30').Value = '  -3.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '172.99'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0836'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.19'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.48%  '
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('E36').Value = '  -5.42%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0343'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.46%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.27'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.27'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.49%  '
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.73'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +17.78%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.24'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '59.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('E44').Value = '  -4.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.39'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.30%  '
$ws.Range('E46').Value = '  -3.73%  '
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('E48').Value = '  +3.90%  '
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('E50').Value = '  -1.67%  '
$ws.Range('D51').Value = '2.434.64'
$ws.Range('E51').Value = '  -0.29%  '
